$p = $ppt.ActivePresentation

# -------------------------------------------------------------------------
# Change 1 (Slide 4, "Products that support DevOps adoption..." paragraph):
# Split the single italic run
#   "incident management systems, configuration management and collaboration platforms;"
# into three runs:
#   (italic)     "incident management systems, configuration management "
#   (not italic) "and"
#   (italic)     " collaboration platforms;"
# -------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange

$needle1 = "incident management systems, configuration management and collaboration platforms;"
$full4 = $tr4.Text
$idx1 = $full4.IndexOf($needle1)

# Cut the whole (currently fully-italic) sentence out of the text...
$whole1 = $tr4.Characters($idx1 + 1, $needle1.Length)
$whole1.Cut() | Out-Null

# ...and re-insert it right after the preceding " and " run (which is NOT
# italic). InsertAfter merges the new characters into that anchor run, so
# the whole sentence starts out life as plain (non-italic) text, matching
# the format of the run it is glued onto instead of the format it used to
# have.
$anchorStart1 = $idx1 - 5 + 1
$anchor1 = $tr4.Characters($anchorStart1, 5)
$anchor1.InsertAfter($needle1) | Out-Null

# Re-locate the (re-inserted) sentence and split it into the three target
# runs, turning italic back on only for the two flanking pieces. The
# middle "and" piece is left completely untouched, so it keeps the
# "no explicit i attribute" plain formatting it inherited above.
$full4b = $tr4.Text
$idx1b = $full4b.IndexOf($needle1)

$part1a = "incident management systems, configuration management "
$part1b = "and"
$part1c = " collaboration platforms;"

$start1a = $idx1b + 1
$start1b = $start1a + $part1a.Length
$start1c = $start1b + $part1b.Length

$run1a = $tr4.Characters($start1a, $part1a.Length)
$run1b = $tr4.Characters($start1b, $part1b.Length)
$run1c = $tr4.Characters($start1c, $part1c.Length)

$run1a.Font.Italic = $true
$run1c.Font.Italic = $true

# -------------------------------------------------------------------------
# Change 2 (Slide 7, "Teams rely on configuration management..." paragraph):
# Split the single (non-italic) run
#   "Teams rely on configuration management for consistent deployment and hosting environments. "
# into three runs:
#   (not italic) "Teams rely on "
#   (italic)     "configuration management"
#   (not italic) " for consistent deployment and hosting environments. "
# -------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(1)
$tr7 = $shp7.TextFrame.TextRange

$needle2 = "Teams rely on configuration management for consistent deployment and hosting environments. "
$full7 = $tr7.Text
$idx2 = $full7.IndexOf($needle2)

$part2a = "Teams rely on "
$part2b = "configuration management"
$part2c = " for consistent deployment and hosting environments. "

$start2a = $idx2 + 1
$start2b = $start2a + $part2a.Length
$start2c = $start2b + $part2b.Length

$run2a = $tr7.Characters($start2a, $part2a.Length)
$run2b = $tr7.Characters($start2b, $part2b.Length)
$run2c = $tr7.Characters($start2c, $part2c.Length)

$run2b.Font.Italic = $true
